# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# columns on each locale's status sheet for the first (13da42f3...) file entry.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("E2").Value = "2016-03-20 12:39:39"
$ws_zhcn.Range("H2").Value = "2016-03-20 12:39:59"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("E2").Value = "2016-03-20 12:39:42"
$ws_dede.Range("H2").Value = "2016-03-20 12:40:06"
